$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.798.34"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "1.871.45"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'0.7318"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").Value = "'241.01"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'0.3139"
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("D9").Value = "'0.07132"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").Value = "'24.44"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("D11").Value = "'0.08153"
$ws.Range("E11").Value = "  -3.38%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7424"
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.879.09"
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("D14").Value = "'5.347"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").Value = "29.787.55"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "'6.006"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("D18").Value = "'247.94"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("D19").Value = "'13.38"
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("D20").Value = "'0.000007803"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'0.9995"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.118.11"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'7.734"
$ws.Range("E24").Value = "  -3.20%  "
$ws.Range("D25").Value = "'0.1539"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("D26").Value = "'9.207"
$ws.Range("D27").Value = "'163.87"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("D28").Value = "'18.54"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").Value = "'2.018"
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("D30").Value = "'1.448"
$ws.Range("E30").Value = "  -1.99%  "
$ws.Range("D31").Value = "'4.519"
$ws.Range("E31").Value = "  -1.87%  "
$ws.Range("D32").Value = "'1.520"
$ws.Range("E32").Value = "  -0.83%  "
$ws.Range("D33").Value = "'4.182"
$ws.Range("E33").Value = "  -2.35%  "
$ws.Range("D34").Value = "'0.05307"
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("D35").Value = "'1.231"
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("D36").Value = "'0.7400"
$ws.Range("E36").Value = "  -2.66%  "
$ws.Range("D37").Value = "'0.9981"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").Value = "'2.702"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").Value = "'0.01936"
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("D40").Value = "'2.731"
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("D41").Value = "'0.4470"
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.8693"
$ws.Range("E42").Value = "  +0.68%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.957"
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("D44").Value = "'71.24"
$ws.Range("E44").Value = "  -1.98%  "
$ws.Range("D45").Value = "1.046.25"
$ws.Range("E45").Value = "  -6.39%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").Value = "'103.81"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("E48").Value = "  -1.56%  "
$ws.Range("D49").Value = "'7.431"
$ws.Range("E49").Value = "  -3.58%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.546"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.018.22"
$ws.Range("E51").Value = "  +0.07%  "
